$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Sheet ALC row 29 (hunk 0)
$ws_ALC.Cells.Item(29, 8).Value = 1100.8572
$ws_ALC.Cells.Item(29, 9).Value = 916.6667
$ws_ALC.Cells.Item(29, 10).Value = 1239
$ws_ALC.Cells.Item(29, 11).Value = 2750.0001
$ws_ALC.Cells.Item(29, 12).Value = 3717
$ws_ALC.Cells.Item(29, 13).Value = -2469.0001
$ws_ALC.Cells.Item(29, 14).Value = -4279

# Sheet ALC row 98 (hunk 1)
$ws_ALC.Cells.Item(98, 8).Value = 7458.8887
$ws_ALC.Cells.Item(98, 9).Value = 6463.7617
$ws_ALC.Cells.Item(98, 10).Value = 8329.625
$ws_ALC.Cells.Item(98, 11).Value = 6463.7617
$ws_ALC.Cells.Item(98, 12).Value = 8329.625
$ws_ALC.Cells.Item(98, 13).Value = -4965.7617
$ws_ALC.Cells.Item(98, 14).Value = -11325.625

# Sheet ALC row 121 (hunk 2)
$ws_ALC.Cells.Item(121, 8).Value = 1500
$ws_ALC.Cells.Item(121, 10).Value = 1500
$ws_ALC.Cells.Item(121, 12).Value = 4500
$ws_ALC.Cells.Item(121, 14).Value = -7994

# Sheet ALC row 122 (hunk 3)
$ws_ALC.Cells.Item(122, 8).Value = 7458.8887
$ws_ALC.Cells.Item(122, 9).Value = 6463.7617
$ws_ALC.Cells.Item(122, 10).Value = 8329.625
$ws_ALC.Cells.Item(122, 11).Value = 19391.2851
$ws_ALC.Cells.Item(122, 12).Value = 24988.875
$ws_ALC.Cells.Item(122, 13).Value = -16941.2851
$ws_ALC.Cells.Item(122, 14).Value = -29888.875

# Sheet ALC row 135 (hunk 4)
$ws_ALC.Cells.Item(135, 8).Value = 673.5
$ws_ALC.Cells.Item(135, 9).Value = 440
$ws_ALC.Cells.Item(135, 10).Value = 1374
$ws_ALC.Cells.Item(135, 11).Value = 3960
$ws_ALC.Cells.Item(135, 12).Value = 12366
$ws_ALC.Cells.Item(135, 13).Value = -1425
$ws_ALC.Cells.Item(135, 14).Value = -17436

# Sheet ALC row 137 (hunk 5)
$ws_ALC.Cells.Item(137, 8).Value = 1163354.2
$ws_ALC.Cells.Item(137, 9).Value = 1702084
$ws_ALC.Cells.Item(137, 11).Value = 5106252
$ws_ALC.Cells.Item(137, 13).Value = -5103702

# Sheet ALC row 138 (hunk 6)
$ws_ALC.Cells.Item(138, 8).Value = 2778.8333
$ws_ALC.Cells.Item(138, 9).Value = 1547.4667
$ws_ALC.Cells.Item(138, 10).Value = 3025.1067
$ws_ALC.Cells.Item(138, 11).Value = 4642.4001
$ws_ALC.Cells.Item(138, 12).Value = 9075.320099999999
$ws_ALC.Cells.Item(138, 13).Value = 497.5999000000002
$ws_ALC.Cells.Item(138, 14).Value = -19355.3201

# Sheet ARM row 45 (hunk 7)
$ws_ARM.Cells.Item(45, 8).Value = 4464.6665
$ws_ARM.Cells.Item(45, 10).Value = 1914
$ws_ARM.Cells.Item(45, 12).Value = 1914
$ws_ARM.Cells.Item(45, 14).Value = -2668

# Sheet ARM row 61 (hunk 8)
$ws_ARM.Cells.Item(61, 8).Value = 1840
$ws_ARM.Cells.Item(61, 9).Value = 1800
$ws_ARM.Cells.Item(61, 11).Value = 1800
$ws_ARM.Cells.Item(61, 13).Value = -1588

# Sheet ARM row 110 (hunk 9)
$ws_ARM.Cells.Item(110, 8).Value = 522.5333000000001
$ws_ARM.Cells.Item(110, 9).Value = 382.30768
$ws_ARM.Cells.Item(110, 10).Value = 1434
$ws_ARM.Cells.Item(110, 11).Value = 382.30768
$ws_ARM.Cells.Item(110, 12).Value = 1434
$ws_ARM.Cells.Item(110, 13).Value = 1662.69232
$ws_ARM.Cells.Item(110, 14).Value = -5524

# Sheet ARM row 136 (hunk 10)
$ws_ARM.Cells.Item(136, 8).Value = 1840
$ws_ARM.Cells.Item(136, 9).Value = 1800
$ws_ARM.Cells.Item(136, 11).Value = 5400
$ws_ARM.Cells.Item(136, 13).Value = -2850

# Sheet BSM row 107 (hunk 11)
$ws_BSM.Cells.Item(107, 8).Value = 1362.3334
$ws_BSM.Cells.Item(107, 9).Value = 1032.4615
$ws_BSM.Cells.Item(107, 10).Value = 2220
$ws_BSM.Cells.Item(107, 11).Value = 1032.4615
$ws_BSM.Cells.Item(107, 12).Value = 2220
$ws_BSM.Cells.Item(107, 13).Value = 887.5385000000001
$ws_BSM.Cells.Item(107, 14).Value = -6060

# Sheet CRP row 31 (hunk 12)
$ws_CRP.Cells.Item(31, 8).Value = 4673.593
$ws_CRP.Cells.Item(31, 10).Value = 12997.857
$ws_CRP.Cells.Item(31, 12).Value = 12997.857
$ws_CRP.Cells.Item(31, 14).Value = -13587.857

# Sheet CRP row 34 (hunk 13)
$ws_CRP.Cells.Item(34, 8).Value = 4673.593
$ws_CRP.Cells.Item(34, 10).Value = 12997.857
$ws_CRP.Cells.Item(34, 12).Value = 12997.857
$ws_CRP.Cells.Item(34, 14).Value = -13401.857

# Sheet CRP row 125 (hunk 14)
$ws_CRP.Cells.Item(125, 8).Value = 35216.668
$ws_CRP.Cells.Item(125, 10).Value = 35216.668
$ws_CRP.Cells.Item(125, 12).Value = 35216.668
$ws_CRP.Cells.Item(125, 14).Value = -40136.668

# Sheet CRP row 137 (hunk 15)
$ws_CRP.Cells.Item(137, 8).Value = 49780
$ws_CRP.Cells.Item(137, 10).Value = 49780
$ws_CRP.Cells.Item(137, 12).Value = 49780
$ws_CRP.Cells.Item(137, 14).Value = -59980

# Sheet CUL row 5 (hunk 16)
$ws_CUL.Cells.Item(5, 8).Value = 835898.75
$ws_CUL.Cells.Item(5, 9).Value = 493.66666
$ws_CUL.Cells.Item(5, 11).Value = 1480.99998
$ws_CUL.Cells.Item(5, 13).Value = -1368.99998

# Sheet CUL row 112 (hunk 17)
$ws_CUL.Cells.Item(112, 8).Value = 1409
$ws_CUL.Cells.Item(112, 10).Value = 2800
$ws_CUL.Cells.Item(112, 12).Value = 8400
$ws_CUL.Cells.Item(112, 14).Value = -10616

# Sheet CUL row 131 (hunk 18)
$ws_CUL.Cells.Item(131, 8).Value = 5155522.5
$ws_CUL.Cells.Item(131, 9).Value = 45455040
$ws_CUL.Cells.Item(131, 10).Value = 932.77905
$ws_CUL.Cells.Item(131, 11).Value = 136365120
$ws_CUL.Cells.Item(131, 12).Value = 2798.33715
$ws_CUL.Cells.Item(131, 13).Value = -136360080
$ws_CUL.Cells.Item(131, 14).Value = -12878.33715

# Sheet CUL row 135 (hunk 19)
$ws_CUL.Cells.Item(135, 8).Value = 835898.75
$ws_CUL.Cells.Item(135, 9).Value = 493.66666
$ws_CUL.Cells.Item(135, 11).Value = 4442.99994
$ws_CUL.Cells.Item(135, 13).Value = -1907.99994

# Sheet CUL row 140 (hunk 20)
$ws_CUL.Cells.Item(140, 8).Value = 3200
$ws_CUL.Cells.Item(140, 9).Value = 3346.1538
$ws_CUL.Cells.Item(140, 11).Value = 10038.4614
$ws_CUL.Cells.Item(140, 13).Value = -4858.4614

# Sheet GSM row 12 (hunk 21)
$ws_GSM.Cells.Item(12, 8).Value = 21400
$ws_GSM.Cells.Item(12, 9).Value = 13000
$ws_GSM.Cells.Item(12, 10).Value = 29800
$ws_GSM.Cells.Item(12, 11).Value = 13000
$ws_GSM.Cells.Item(12, 12).Value = 29800
$ws_GSM.Cells.Item(12, 13).Value = -12860
$ws_GSM.Cells.Item(12, 14).Value = -30080

# Sheet GSM row 80 (hunk 22)
$ws_GSM.Cells.Item(80, 8).Value = 27782322
$ws_GSM.Cells.Item(80, 9).Value = 50004700
$ws_GSM.Cells.Item(80, 10).Value = 4350
$ws_GSM.Cells.Item(80, 11).Value = 50004700
$ws_GSM.Cells.Item(80, 12).Value = 4350
$ws_GSM.Cells.Item(80, 13).Value = -50003702
$ws_GSM.Cells.Item(80, 14).Value = -6346

# Sheet GSM row 83 (hunk 23)
$ws_GSM.Cells.Item(83, 8).Value = 27782322
$ws_GSM.Cells.Item(83, 9).Value = 50004700
$ws_GSM.Cells.Item(83, 10).Value = 4350
$ws_GSM.Cells.Item(83, 11).Value = 250023500
$ws_GSM.Cells.Item(83, 12).Value = 21750
$ws_GSM.Cells.Item(83, 13).Value = -250018508
$ws_GSM.Cells.Item(83, 14).Value = -31734

# Sheet LTW row 55 (hunk 24)
$ws_LTW.Cells.Item(55, 8).Value = 365.91666
$ws_LTW.Cells.Item(55, 9).Value = 112.888885
$ws_LTW.Cells.Item(55, 10).Value = 1125
$ws_LTW.Cells.Item(55, 11).Value = 112.888885
$ws_LTW.Cells.Item(55, 12).Value = 1125
$ws_LTW.Cells.Item(55, 13).Value = 60.111115
$ws_LTW.Cells.Item(55, 14).Value = -1471

# Sheet LTW row 82 (hunk 25)
$ws_LTW.Cells.Item(82, 8).Value = 1324.5883
$ws_LTW.Cells.Item(82, 9).Value = 814.24
$ws_LTW.Cells.Item(82, 10).Value = 2742.2222
$ws_LTW.Cells.Item(82, 11).Value = 814.24
$ws_LTW.Cells.Item(82, 12).Value = 2742.2222
$ws_LTW.Cells.Item(82, 13).Value = -453.24
$ws_LTW.Cells.Item(82, 14).Value = -3464.2222

# Sheet LTW row 85 (hunk 26)
$ws_LTW.Cells.Item(85, 8).Value = 1324.5883
$ws_LTW.Cells.Item(85, 9).Value = 814.24
$ws_LTW.Cells.Item(85, 10).Value = 2742.2222
$ws_LTW.Cells.Item(85, 11).Value = 814.24
$ws_LTW.Cells.Item(85, 12).Value = 2742.2222
$ws_LTW.Cells.Item(85, 13).Value = 433.76
$ws_LTW.Cells.Item(85, 14).Value = -5238.2222

# Sheet LTW row 125 (hunk 27)
$ws_LTW.Cells.Item(125, 8).Value = 41805.625
$ws_LTW.Cells.Item(125, 10).Value = 41805.625
$ws_LTW.Cells.Item(125, 12).Value = 41805.625
$ws_LTW.Cells.Item(125, 14).Value = -51645.625

# Sheet LTW row 136 (hunk 28)
$ws_LTW.Cells.Item(136, 8).Value = 4264.316
$ws_LTW.Cells.Item(136, 9).Value = 1474.7273
$ws_LTW.Cells.Item(136, 11).Value = 4424.1819
$ws_LTW.Cells.Item(136, 13).Value = -1874.1819

# Sheet WVR row 4 (hunk 29)
$ws_WVR.Cells.Item(4, 8).Value = 3908.4285
$ws_WVR.Cells.Item(4, 9).Value = 700
$ws_WVR.Cells.Item(4, 11).Value = 700
$ws_WVR.Cells.Item(4, 13).Value = -587

# Sheet WVR row 9 (hunk 30)
$ws_WVR.Cells.Item(9, 8).Value = 3000
$ws_WVR.Cells.Item(9, 10).Value = 0
$ws_WVR.Cells.Item(9, 12).Value = 0
$ws_WVR.Cells.Item(9, 14).ClearContents()

# Sheet WVR row 62 (hunk 31)
$ws_WVR.Cells.Item(62, 8).Value = 11070.714
$ws_WVR.Cells.Item(62, 9).Value = 3810
$ws_WVR.Cells.Item(62, 10).Value = 24140
$ws_WVR.Cells.Item(62, 11).Value = 3810
$ws_WVR.Cells.Item(62, 12).Value = 24140
$ws_WVR.Cells.Item(62, 13).Value = -3186
$ws_WVR.Cells.Item(62, 14).Value = -25388

# Sheet WVR row 65 (hunk 32)
$ws_WVR.Cells.Item(65, 8).Value = 11070.714
$ws_WVR.Cells.Item(65, 9).Value = 3810
$ws_WVR.Cells.Item(65, 10).Value = 24140
$ws_WVR.Cells.Item(65, 11).Value = 19050
$ws_WVR.Cells.Item(65, 12).Value = 120700
$ws_WVR.Cells.Item(65, 13).Value = -15930
